# Add GA4 start dates to both sheets (NAV + 5525)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("NAV")
$ws2 = $wb.Worksheets.Item("5525")

# NAV sheet: row 7 -> ga4_start / 2025-02-25
$ws1.Range("A7").Value = "ga4_start"
$ws1.Range("B7").Value = 45713
$ws1.Range("B7").NumberFormat = "mm-dd-yy"

# 5525 sheet: row 7 -> ga4_start / 2025-08-01
$ws2.Range("A7").Value = "ga4_start"
$ws2.Range("B7").Value = 45870

# Re-use the exact same date style that was just created on the NAV sheet
# (instead of setting NumberFormat again, which would mint a second,
# duplicate cellXfs entry) by copying the formatting over.
$ws1.Range("B7").Copy()
$ws2.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update selections to the newly added cell on each sheet, and make sure
# the "5525" sheet is touched/selected first so that the final active tab
# (and the one left with tabSelected="1") ends up being "NAV", matching
# the target workbook state.
$ws2.Range("B7").Select()
$ws1.Range("B7").Select()
